$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price / 1h-volume figures (incl. the Kaspa <-> Cosmos
# row swap at rows 41-42). Cells that look like plain numbers get a leading
# apostrophe so Excel keeps them as text, matching the source data which
# stores every Price/Volume cell as a string (e.g. '62.004.40').

$ws.Range("D2").Value = "61.956.93"
$ws.Range("E2").Value = "  -1.12%  "
$ws.Range("D3").Value = "2.908.87"
$ws.Range("E3").Value = "  -1.92%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'579.23"
$ws.Range("E5").Value = "  -2.53%  "
$ws.Range("D6").Value = "'146.16"
$ws.Range("E6").Value = "  +0.57%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").Value = "'0.505"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "2.908.53"
$ws.Range("E9").Value = "  -1.92%  "
$ws.Range("D10").Value = "'6.71"
$ws.Range("E10").Value = "  -7.07%  "
$ws.Range("E11").Value = "  +4.36%  "
$ws.Range("E12").Value = "  -2.58%  "
$ws.Range("D13").Value = "'0.0000237"
$ws.Range("E13").Value = "  -0.17%  "
$ws.Range("D14").Value = "'32.72"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("E15").Value = "  -0.81%  "
$ws.Range("D16").Value = "3.393.60"
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("D17").Value = "61.960.85"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").Value = "'6.65"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").Value = "2.908.18"
$ws.Range("E19").Value = "  -2.11%  "
$ws.Range("D20").Value = "'435.46"
$ws.Range("E20").Value = "  -1.33%  "
$ws.Range("E21").Value = "  -0.96%  "
$ws.Range("D22").Value = "'0.660"
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("D23").Value = "'6.96"
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("D24").Value = "'79.95"
$ws.Range("E24").Value = "  -2.07%  "
$ws.Range("D26").Value = "'10.19"
$ws.Range("E26").Value = "  -9.35%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'2.05"
$ws.Range("E28").Value = "  -3.42%  "
$ws.Range("D29").Value = "'0.0000110"
$ws.Range("E29").Value = "  +16.10%  "
$ws.Range("D30").Value = "'7.08"
$ws.Range("E30").Value = "  -2.15%  "
$ws.Range("E31").Value = "  -2.33%  "
$ws.Range("D32").Value = "'2.10"
$ws.Range("E32").Value = "  -1.59%  "
$ws.Range("D33").Value = "'0.107"
$ws.Range("E33").Value = "  -2.53%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'25.80"
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("D36").Value = "'0.969"
$ws.Range("E36").Value = "  -2.80%  "
$ws.Range("D37").Value = "'3.08"
$ws.Range("E37").Value = "  +1.98%  "
$ws.Range("D38").Value = "'5.49"
$ws.Range("E38").Value = "  -2.39%  "
$ws.Range("D39").Value = "'49.20"
$ws.Range("E39").Value = "  -0.48%  "
$ws.Range("E40").Value = "  -3.02%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").Value = "'0.116"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "'8.31"
$ws.Range("E42").Value = "  -2.69%  "
$ws.Range("E43").Value = "  -4.10%  "
$ws.Range("D44").Value = "'38.50"
$ws.Range("E44").Value = "  -3.86%  "
$ws.Range("D45").Value = "2.688.44"
$ws.Range("E45").Value = "  -2.01%  "
$ws.Range("D46").Value = "'134.87"
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").Value = "'343.86"
$ws.Range("E48").Value = "  -4.96%  "
$ws.Range("E50").Value = "  -1.67%  "
$ws.Range("D51").Value = "'21.92"
$ws.Range("E51").Value = "  -4.74%  "
